$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Row = 7; I = "sv"; J = "Statement-opinion" }
    @{ Row = 10; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 11; I = "ba"; J = "Appreciation" }
    @{ Row = 19; I = "sv"; J = "Statement-opinion" }
    @{ Row = 21; I = "aa"; J = "Agree/Accept" }
    @{ Row = 32; I = "sv"; J = "Statement-opinion" }
    @{ Row = 35; I = "aa"; J = "Agree/Accept" }
    @{ Row = 48; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 90; I = "%"; J = "Uninterpretable" }
    @{ Row = 94; I = "ba"; J = "Appreciation" }
    @{ Row = 107; I = "b"; J = "Acknowledge (Backchannel)" }
    @{ Row = 111; I = "aa"; J = "Agree/Accept" }
    @{ Row = 112; I = "sv"; J = "Statement-opinion" }
    @{ Row = 114; I = "sv"; J = "Statement-opinion" }
    @{ Row = 117; I = "sv"; J = "Statement-opinion" }
    @{ Row = 121; I = "%"; J = "Uninterpretable" }
    @{ Row = 122; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 133; I = "sv"; J = "Statement-opinion" }
    @{ Row = 137; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 139; I = "aa"; J = "Agree/Accept" }
    @{ Row = 159; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 168; I = "sv"; J = "Statement-opinion" }
    @{ Row = 171; I = "sv"; J = "Statement-opinion" }
    @{ Row = 172; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 181; I = "sv"; J = "Statement-opinion" }
    @{ Row = 201; I = "sv"; J = "Statement-opinion" }
    @{ Row = 205; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 232; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 237; I = "aa"; J = "Agree/Accept" }
    @{ Row = 260; I = "%"; J = "Uninterpretable" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}

